$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.889.17'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.02%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.487.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.85%  '

# Row 4
$ws.Range('E4').Value = '  +0.23%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.03%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.15%  '

# Row 7
$ws.Range('E7').Value = '  -0.20%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.540'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.48%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.115'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.88%  '

# Row 10
$ws.Range('E10').Value = '  +0.45%  '

# Row 11
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.368'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.62%  '

# Row 12
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.38'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.99%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.38%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000187'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.33%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.935.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.75%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.372.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.24%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.543.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.06%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.52%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.17%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '330.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.11%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.45%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.08%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.01%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.91%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '658.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +16.30%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.27%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000107'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.35%  '

# Row 28
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.69%  '

# Row 29
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.601.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.64%  '

# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.64'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.99%  '

# Row 31
$ws.Range('E31').Value = '  -0.51%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.147'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.86%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.17%  '

# Row 34
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.26'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.70%  '

# Row 35
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.74%  '

# Row 36
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.390'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.36%  '

# Row 37
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.997'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.24%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.64'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.12'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.88%  '

# Row 40
$ws.Range('E40').Value = '  +2.46%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '148.67'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.78%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.21%  '

# Row 43
$ws.Range('E43').Value = '  +0.03%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '154.27'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.65%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.50%  '

# Row 46
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.56%  '

# Row 47
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0553'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.34%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.615'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.43%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0240'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.16%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0932'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.11%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.36%  '
